$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 18521136
$ws.Range("I76").Value = 111111110
$ws.Range("J76").Value = 3140
$ws.Range("K76").Value = 111111110
$ws.Range("L76").Value = 3140
$ws.Range("M76").Value = -111110795
$ws.Range("N76").Value = -3770

$ws.Range("H79").Value = 18521136
$ws.Range("I79").Value = 111111110
$ws.Range("J79").Value = 3140
$ws.Range("K79").Value = 111111110
$ws.Range("L79").Value = 3140
$ws.Range("M79").Value = -111110018
$ws.Range("N79").Value = -5324

$ws.Range("H80").Value = 605.45
$ws.Range("I80").Value = 300.1111
$ws.Range("J80").Value = 855.2727
$ws.Range("K80").Value = 900.3333
$ws.Range("L80").Value = 2565.8181
$ws.Range("M80").Value = 97.6667
$ws.Range("N80").Value = -4561.8181

$ws.Range("H83").Value = 605.45
$ws.Range("I83").Value = 300.1111
$ws.Range("J83").Value = 855.2727
$ws.Range("K83").Value = 2700.9999
$ws.Range("L83").Value = 7697.4543
$ws.Range("M83").Value = 2291.0001
$ws.Range("N83").Value = -17681.4543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6892.6
$ws.Range("I2").Value = 9369.214
$ws.Range("J2").Value = 1113.8334
$ws.Range("K2").Value = 9369.214
$ws.Range("L2").Value = 1113.8334
$ws.Range("M2").Value = -9256.214
$ws.Range("N2").Value = -1339.8334

$ws.Range("H45").Value = 2262.4119
$ws.Range("I45").Value = 1232.8182
$ws.Range("J45").Value = 4150
$ws.Range("K45").Value = 1232.8182
$ws.Range("L45").Value = 4150
$ws.Range("M45").Value = -855.8182
$ws.Range("N45").Value = -4904

$ws.Range("H88").Value = 5370
$ws.Range("I88").Value = 7616.6665
$ws.Range("K88").Value = 7616.6665
$ws.Range("M88").Value = -7210.6665

$ws.Range("H91").Value = 5370
$ws.Range("I91").Value = 7616.6665
$ws.Range("K91").Value = 7616.6665
$ws.Range("M91").Value = -6212.6665

$ws.Range("H116").Value = 6892.6
$ws.Range("I116").Value = 9369.214
$ws.Range("J116").Value = 1113.8334
$ws.Range("K116").Value = 9369.214
$ws.Range("L116").Value = 1113.8334
$ws.Range("M116").Value = -7075.214
$ws.Range("N116").Value = -5701.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6892.6
$ws.Range("I3").Value = 9369.214
$ws.Range("J3").Value = 1113.8334
$ws.Range("K3").Value = 9369.214
$ws.Range("L3").Value = 1113.8334
$ws.Range("M3").Value = -9255.214
$ws.Range("N3").Value = -1341.8334

$ws.Range("H5").Value = 525000
$ws.Range("I5").Value = 525000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 525000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -524887
$ws.Range("N5").ClearContents()

$ws.Range("H86").Value = 19406.889
$ws.Range("I86").Value = 10165
$ws.Range("J86").Value = 51753.5
$ws.Range("K86").Value = 10165
$ws.Range("L86").Value = 51753.5
$ws.Range("M86").Value = -9042
$ws.Range("N86").Value = -53999.5

$ws.Range("H89").Value = 19406.889
$ws.Range("I89").Value = 10165
$ws.Range("J89").Value = 51753.5
$ws.Range("K89").Value = 50825
$ws.Range("L89").Value = 258767.5
$ws.Range("M89").Value = -45209
$ws.Range("N89").Value = -269999.5

$ws.Range("H105").Value = 3130.4348
$ws.Range("I105").Value = 2823.5293
$ws.Range("K105").Value = 2823.5293
$ws.Range("M105").Value = -1076.5293

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 126250560
$ws.Range("I4").Value = 500.33334
$ws.Range("J4").Value = 202000600
$ws.Range("K4").Value = 500.33334
$ws.Range("L4").Value = 202000600
$ws.Range("M4").Value = -388.33334
$ws.Range("N4").Value = -202000824

$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498

$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488

$ws.Range("H107").Value = 320.875
$ws.Range("I107").Value = 250.5
$ws.Range("J107").Value = 391.25
$ws.Range("K107").Value = 250.5
$ws.Range("L107").Value = 391.25
$ws.Range("M107").Value = 1669.5
$ws.Range("N107").Value = -4231.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 224.5
$ws.Range("I4").Value = 224.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 673.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -561.5
$ws.Range("N4").ClearContents()

$ws.Range("H129").Value = 1043.25
$ws.Range("I129").Value = 315
$ws.Range("J129").Value = 1188.9
$ws.Range("K129").Value = 945
$ws.Range("L129").Value = 3566.7
$ws.Range("M129").Value = 4055
$ws.Range("N129").Value = -13566.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H97").Value = 2408.1667
$ws.Range("I97").Value = 2408.1667
$ws.Range("K97").Value = 2408.1667
$ws.Range("M97").Value = -1912.1667

$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280

$ws.Range("H139").Value = 59519.6
$ws.Range("J139").Value = 59519.6
$ws.Range("L139").Value = 59519.6
$ws.Range("N139").Value = -69799.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 46816.184
$ws.Range("I82").Value = 72642.57
$ws.Range("J82").Value = 1620
$ws.Range("K82").Value = 72642.57
$ws.Range("L82").Value = 1620
$ws.Range("M82").Value = -72281.57
$ws.Range("N82").Value = -2342

$ws.Range("H85").Value = 46816.184
$ws.Range("I85").Value = 72642.57
$ws.Range("J85").Value = 1620
$ws.Range("K85").Value = 72642.57
$ws.Range("L85").Value = 1620
$ws.Range("M85").Value = -71394.57
$ws.Range("N85").Value = -4116

$ws.Range("H132").Value = 4111.9414
$ws.Range("I132").Value = 3878.1365
$ws.Range("J132").Value = 4540.5835
$ws.Range("K132").Value = 11634.4095
$ws.Range("L132").Value = 13621.7505
$ws.Range("M132").Value = -9104.4095
$ws.Range("N132").Value = -18681.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 36387410
$ws.Range("I2").Value = 50006436
$ws.Range("J2").Value = 69998.664
$ws.Range("K2").Value = 50006436
$ws.Range("L2").Value = 69998.664
$ws.Range("M2").Value = -50006324
$ws.Range("N2").Value = -70222.664
